$wb = $excel.ActiveWorkbook
$players = $wb.Worksheets.Item("Players")
$ownerTotals = $wb.Worksheets.Item("OwnerTotals")

# Set column G width so it renders as 17 characters wide after COM pixel rounding
$players.Columns.Item(7).ColumnWidth = 16.17

$players.Range("G4").Value = '8:15 - 2nd Half'
$players.Range("H4").Value = 12
$players.Range("O4").Value = 30
$players.Range("G9").Value = '8:15 - 2nd Half'
$players.Range("H9").Value = 2
$players.Range("I9").Value = 6
$players.Range("O9").Value = 20
$players.Range("D11").Value = 'Nijel Pack'
$players.Range("E11").Value = 'OU'
$players.Range("F11").Value = 'OU@MIZ'
$players.Range("G11").Value = '8:15 - 2nd Half'
$players.Range("H11").Value = 14
$players.Range("I11").Value = 13
$players.Range("J11").Value = 3
$players.Range("K11").Value = 5
$players.Range("N11").Value = 3
$players.Range("O11").Value = 24
$players.Range("D12").Value = 'Malik Dia'
$players.Range("E12").Value = 'MISS'
$players.Range("F12").Value = 'MISS@UK'
$players.Range("G12").Value = 'Final'
$players.Range("I12").Value = 16
$players.Range("J12").Value = 4
$players.Range("K12").Value = 0
$players.Range("N12").Value = 0
$players.Range("O12").Value = 21
$players.Range("G15").Value = '8:15 - 2nd Half'
$players.Range("H15").Value = 22
$players.Range("I15").Value = 17
$players.Range("K15").Value = 3
$players.Range("O15").Value = 30
$players.Range("G18").Value = '8:15 - 2nd Half'
$players.Range("H18").Value = 10
$players.Range("I18").Value = 7
$players.Range("J18").Value = 4
$players.Range("O18").Value = 18
$players.Range("G22").Value = '8:15 - 2nd Half'
$players.Range("H22").Value = 12
$players.Range("I22").Value = 15
$players.Range("J22").Value = 4
$players.Range("O22").Value = 27
$players.Range("G25").Value = '8:15 - 2nd Half'
$players.Range("H25").Value = 10
$players.Range("I25").Value = 10
$players.Range("O25").Value = 16
$players.Range("G26").Value = '8:15 - 2nd Half'
$players.Range("G29").Value = '8:15 - 2nd Half'
$players.Range("H29").Value = 9
$players.Range("N29").Value = 4
$players.Range("O29").Value = 24
$players.Range("G32").Value = '8:15 - 2nd Half'
$players.Range("O32").Value = 12
$players.Range("G40").Value = '8:15 - 2nd Half'
$players.Range("G44").Value = '8:15 - 2nd Half'
$players.Range("H44").Value = 2
$players.Range("O44").Value = 21
$players.Range("D52").Value = 'T.O. Barrett'
$players.Range("E52").Value = 'MIZ'
$players.Range("F52").Value = 'OU@MIZ'
$players.Range("G52").Value = '8:15 - 2nd Half'
$players.Range("H52").Value = 14
$players.Range("I52").Value = 11
$players.Range("J52").Value = 5
$players.Range("K52").Value = 1
$players.Range("L52").Value = 3
$players.Range("M52").Value = 0
$players.Range("O52").Value = 26
$players.Range("D53").Value = 'Andrija Jelavic'
$players.Range("E53").Value = 'UK'
$players.Range("F53").Value = 'MISS@UK'
$players.Range("I53").Value = 6
$players.Range("J53").Value = 7
$players.Range("L53").Value = 1
$players.Range("M53").Value = 1
$players.Range("N53").Value = 2
$players.Range("O53").Value = 22
$players.Range("D54").Value = 'Pop Isaacs'
$players.Range("E54").Value = 'TA&M'
$players.Range("F54").Value = 'SC@TA&M'
$players.Range("G54").Value = 'Final'
$players.Range("H54").Value = 11
$players.Range("I54").Value = 11
$players.Range("J54").Value = 4
$players.Range("K54").Value = 3
$players.Range("L54").Value = 0
$players.Range("N54").Value = 1
$players.Range("O54").Value = 24
$players.Range("D55").Value = 'Shawn Phillips Jr.'
$players.Range("E55").Value = 'MIZ'
$players.Range("F55").Value = 'OU@MIZ'
$players.Range("G55").Value = '8:15 - 2nd Half'
$players.Range("H55").Value = 10
$players.Range("I55").Value = 4
$players.Range("J55").Value = 6
$players.Range("K55").Value = 0
$players.Range("M55").Value = 3
$players.Range("O55").Value = 16
$players.Range("D56").Value = 'Chandler Bing'
$players.Range("E56").Value = 'VAN'
$players.Range("F56").Value = 'VAN@MSST'
$players.Range("I56").Value = 5
$players.Range("J56").Value = 5
$players.Range("K56").Value = 2
$players.Range("L56").Value = 0
$players.Range("M56").Value = 0
$players.Range("N56").Value = 1
$players.Range("O56").Value = 20
$players.Range("D57").Value = 'Chendall Weaver'
$players.Range("E57").Value = 'TEX'
$players.Range("F57").Value = 'UGA@TEX'
$players.Range("I57").Value = 4
$players.Range("J57").Value = 2
$players.Range("K57").Value = 4
$players.Range("L57").Value = 1
$players.Range("O57").Value = 22
$players.Range("D58").Value = 'EJ Walker'
$players.Range("E58").Value = 'SC'
$players.Range("H58").Value = 9
$players.Range("I58").Value = 5
$players.Range("J58").Value = 3
$players.Range("K58").Value = 1
$players.Range("M58").Value = 1
$players.Range("O58").Value = 19
$players.Range("D59").Value = 'Ali Dibba'
$players.Range("E59").Value = 'TA&M'
$players.Range("F59").Value = 'SC@TA&M'
$players.Range("I59").Value = 9
$players.Range("J59").Value = 2
$players.Range("K59").Value = 0
$players.Range("M59").Value = 0
$players.Range("O59").Value = 16
$players.Range("D60").Value = 'Augusto Cassiá'
$players.Range("I60").Value = 4
$players.Range("J60").Value = 1
$players.Range("K60").Value = 2
$players.Range("L60").Value = 0
$players.Range("M60").Value = 1
$players.Range("N60").Value = 0
$players.Range("O60").Value = 10
$players.Range("D61").Value = 'Corey Chest'
$players.Range("E61").Value = 'MISS'
$players.Range("F61").Value = 'MISS@UK'
$players.Range("J61").Value = 10
$players.Range("K61").Value = 0
$players.Range("L61").Value = 1
$players.Range("M61").Value = 0
$players.Range("N61").Value = 2
$players.Range("O61").Value = 21
$players.Range("D62").Value = 'Jamarion Davis-Fleming'
$players.Range("E62").Value = 'MSST'
$players.Range("F62").Value = 'VAN@MSST'
$players.Range("I62").Value = 3
$players.Range("J62").Value = 8
$players.Range("M62").Value = 2
$players.Range("O62").Value = 25
$players.Range("D63").Value = 'Jasper Johnson'
$players.Range("E63").Value = 'UK'
$players.Range("F63").Value = 'MISS@UK'
$players.Range("I63").Value = 11
$players.Range("J63").Value = 2
$players.Range("K63").Value = 2
$players.Range("M63").Value = 0
$players.Range("N63").Value = 3
$players.Range("O63").Value = 17
$players.Range("D64").Value = 'Justin Abson'
$players.Range("E64").Value = 'UGA'
$players.Range("F64").Value = 'UGA@TEX'
$players.Range("G64").Value = 'Final'
$players.Range("I64").Value = 6
$players.Range("J64").Value = 3
$players.Range("M64").Value = 1
$players.Range("O64").Value = 15
$players.Range("G79").Value = '8:15 - 2nd Half'
$players.Range("G86").Value = '8:15 - 2nd Half'
$players.Range("G93").Value = '8:15 - 2nd Half'
$players.Range("G95").Value = '8:15 - 2nd Half'
$players.Range("G98").Value = '8:15 - 2nd Half'
$players.Range("G99").Value = '8:15 - 2nd Half'
$players.Range("D102").Value = 'Ja''Borri McGhee'
$players.Range("E102").Value = 'MSST'
$players.Range("F102").Value = 'VAN@MSST'
$players.Range("G102").Value = 'Final'
$players.Range("H102").Value = -4
$players.Range("I102").Value = 2
$players.Range("J102").Value = 0
$players.Range("K102").Value = 1
$players.Range("O102").Value = 17
$players.Range("D103").Value = 'Trent Pierce'
$players.Range("E103").Value = 'MIZ'
$players.Range("F103").Value = 'OU@MIZ'
$players.Range("G103").Value = '8:15 - 2nd Half'
$players.Range("I103").Value = 0
$players.Range("J103").Value = 1
$players.Range("K103").Value = 0
$players.Range("O103").Value = 7

$ownerTotals.Range("B2").Value = 55
$ownerTotals.Range("B3").Value = 55
$ownerTotals.Range("B5").Value = 43
$ownerTotals.Range("B6").Value = 35
$ownerTotals.Range("B8").Value = 25
